# Refresh the cryptos table's Price (D) and Volume(1h) (E) columns
# per the Tue Sep 24 13:43:54 UTC 2024 GitHub Actions data refresh.
# Source data is plain text (e.g. "0.999", "  +2.21%  "); some of the new
# Price strings are numeric-looking, so those cells are forced to text
# first (NumberFormat "@") and the style is reset back to Normal afterwards
# so only the displayed text -- not the cell formatting -- changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '63.540.28'
    'E2' = '  +0.45%  '
    'D3' = '2.637.04'
    'E3' = '  -0.12%  '
    'D4' = '0.999'
    'E4' = '  -0.08%  '
    'D5' = '607.20'
    'E5' = '  +2.21%  '
    'D6' = '147.56'
    'E6' = '  +2.86%  '
    'D7' = '0.999'
    'E7' = '  -0.08%  '
    'D8' = '0.590'
    'E9' = '  +2.41%  '
    'E10' = '  -1.46%  '
    'D11' = '0.373'
    'E11' = '  +5.13%  '
    'E12' = '  -0.12%  '
    'D13' = '27.57'
    'E13' = '  +0.79%  '
    'D14' = '3.106.25'
    'E14' = '  -0.36%  '
    'D15' = '63.365.48'
    'E15' = '  +0.33%  '
    'E16' = '  +2.79%  '
    'D17' = '2.640.61'
    'E17' = '  +0.38%  '
    'D18' = '11.59'
    'E18' = '  +1.91%  '
    'D19' = '4.58'
    'E19' = '  +5.09%  '
    'D20' = '345.21'
    'E20' = '  +1.84%  '
    'E21' = '  +2.59%  '
    'D22' = '0.999'
    'E22' = '  -0.17%  '
    'D23' = '5.58'
    'E23' = '  -3.43%  '
    'D24' = '66.93'
    'E24' = '  +0.13%  '
    'E25' = '  +2.58%  '
    'D26' = '9.09'
    'E26' = '  +7.92%  '
    'D27' = '1.58'
    'E27' = '  +3.02%  '
    'D28' = '562.91'
    'E28' = '  +7.30%  '
    'D29' = '8.07'
    'E29' = '  +3.91%  '
    'E30' = '  -0.98%  '
    'D31' = '1.00'
    'E31' = '  -0.18%  '
    'D32' = '2.06'
    'E32' = '  +4.40%  '
    'D33' = '0.0₃0856'
    'E33' = '  +6.25%  '
    'E34' = '  -2.79%  '
    'D35' = '5.19'
    'E35' = '  +5.94%  '
    'D36' = '167.48'
    'E36' = '  -3.94%  '
    'D37' = '0.407'
    'E37' = '  +1.21%  '
    'E38' = '  -0.11%  '
    'E39' = '  +8.93%  '
    'D40' = '19.17'
    'E40' = '  +0.88%  '
    'E41' = '  +0.06%  '
    'D42' = '166.04'
    'E42' = '  -3.13%  '
    'E43' = '  +1.94%  '
    'D44' = '22.18'
    'E44' = '  +1.25%  '
    'D45' = '0.0573'
    'E45' = '  +3.03%  '
    'D46' = '0.631'
    'E46' = '  +0.28%  '
    'E47' = '  +3.94%  '
    'E48' = '  +0.38%  '
    'D49' = '1.94'
    'E49' = '  +14.14%  '
    'D50' = '18.96'
    'E50' = '  +2.67%  '
    'D51' = '0.184'
    'E51' = '  +6.16%  '
}

# Cells whose new value would otherwise be auto-converted to a number by Excel.
$textCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D11', 'D13', 'D18', 'D19', 'D20', 'D22', 'D23', 'D24', 'D26', 'D27', 'D28', 'D29', 'D31', 'D32', 'D35', 'D36', 'D37', 'D40', 'D42', 'D44', 'D45', 'D46', 'D49', 'D50', 'D51')

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    if ($textCells -contains $addr) {
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$addr]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $updates[$addr]
    }
}
